$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.755.14'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '1.904.55'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9970'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.48'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9984'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5235'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +6.50%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3780'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07220'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.28'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9077'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07628'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.909.29'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.442'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.05'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9968'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008682'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9978'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("D19").Value = '27.784.91'
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.143'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").Value = '2.166.92'
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.84'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.605'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("E26").Value = '  -2.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.165'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.29'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("E29").Value = '  -1.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.837'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08998'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.880'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +5.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.175'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7771'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02090'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.72%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.618'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.064'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.092'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5521'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05271'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.669'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '114.62'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.536'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1512'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4800'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.45'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9982'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("E49").Value = '  -0.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '66.78'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05991'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.99%  '
